# Modified DSL for PB tests
# Applies the content edits to the "TestCases" sheet of EBPB_Imager.xlsx:
#  - Row 2 ("Change Start page"): Subject renamed from "Signature-JS-4.1" to "PB-Imager";
#    Results cleared (was "Pass").
#  - Rows 3-6: Results cleared (was "Pass").
#  - Rows 9-15: Testcase ID (column A) renumbered sequentially (8, 9, 10, 11, 12, 13, 14).
#  - Rows 10, 11, 12, 13, 14, 15 (Steps column G): wait(7) -> wait(15).
#  - Row 12 (Expected Behaviour column H, VT056-1140): validate_Result filled in with
#    "Error Code: 12014" (was empty).
#  - Row 12 (column A) also loses its stray cell-shading, matching the style already used
#    by the other rows in that column.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestCases")

# --- Row 2: Subject + Result ---
$ws.Range("C2").Value = "PB-Imager"
$ws.Range("J2").Value = ""

# --- Rows 3-6: clear stale "Pass" results ---
$ws.Range("J3").Value = ""
$ws.Range("J4").Value = ""
$ws.Range("J5").Value = ""
$ws.Range("J6").Value = ""

# --- Row 9: renumber Testcase ID ---
$ws.Range("A9").Value = 8

# --- Row 10 (VT056-1138): renumber + wait(7) -> wait(15) ---
$ws.Range("A10").Value = 9
$ws.Range("G10").Value = "wait(3);`nvalidate1;`nlink_Click(imager_test_link);`nvalidate2;`nSelectTestToRun(VT056_1138_string);`nvalidate3;`nClickRunTest(runtest_top_xpath);`nwait(15);`nvalidate4;"

# --- Row 11 (VT056-1139): renumber + wait(7) -> wait(15) ---
$ws.Range("A11").Value = 10
$ws.Range("G11").Value = "wait(3);`nvalidate1;`nlink_Click(imager_test_link);`nvalidate2;`nSelectTestToRun(VT056_1139_string);`nvalidate3;`nClickRunTest(runtest_top_xpath);`nwait(15);`nvalidate4;"

# --- Row 12 (VT056-1140): renumber, fix shading, wait(7) -> wait(15), fill in Result ---
$ws.Range("A13").Copy()
$ws.Range("A12").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A12").Value = 11
$ws.Range("G12").Value = "wait(3);`nvalidate1;`nlink_Click(imager_test_link);`nvalidate2;`nwait(5);`nSelectImager(back_camera);`nSelectTestToRun(VT056_1140_string);`nvalidate3;`nClickRunTest(runtest_top_xpath);`nwait(15);`nvalidate4;"
$ws.Range("H12").Value = "validate1`n{`nvalidate_PageTitle=Pocket Browser Tests`n};`nvalidate2`n{`nvalidate_PageTitle=PB and RE2.2 Semi Auto Frame Work : Imager`n};`nvalidate3`n{`nvalidate_Text_Exists=VT056-1140`n};`nvalidate4`n{`nvalidate_Result=Error Code: 12014`n};"

# --- Row 13 (VT056-1142): renumber + wait(7) -> wait(15) ---
$ws.Range("A13").Value = 12
$ws.Range("G13").Value = "wait(3);`nvalidate1;`nlink_Click(imager_test_link);`nvalidate2;`nSelectTestToRun(VT056_1142_string);`nvalidate3;`nClickRunTest(runtest_top_xpath);`nwait(15);`nvalidate4;"

# --- Row 14 (VT056-1143): renumber + wait(7) -> wait(15) ---
$ws.Range("A14").Value = 13
$ws.Range("G14").Value = "wait(3);`nvalidate1;`nlink_Click(imager_test_link);`nvalidate2;`nSelectTestToRun(VT056_1143_string);`nvalidate3;`nClickRunTest(runtest_top_xpath);`nwait(15);`nvalidate4;"

# --- Row 15 (VT056-1146): renumber + wait(7) -> wait(15) ---
$ws.Range("A15").Value = 14
$ws.Range("G15").Value = "wait(3);`nvalidate1;`nlink_Click(imager_test_link);`nvalidate2;`nSelectTestToRun(VT056_1146_string);`nvalidate3;`nClickRunTest(runtest_top_xpath);`nwait(15);`nvalidate4;"
